# DEAN import format workbook update
# - Adds a new Person row (Erica Watley, person_id 10001)
# - Adds an instructor_id column to Course_Section with value 10001
# - Adds a new (mostly blank) Enrollment row with enrollment_id 1
# - Updates sheet view / selection state (tracks the diff's sheetView/selection changes)

$wb = $excel.ActiveWorkbook

# --- Person sheet: append row 3 ---
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Cells.Item(3, 1).Value = 10001
$wsPerson.Cells.Item(3, 4).Value = "Erica"
$wsPerson.Cells.Item(3, 6).Value = "Watley"

# --- Course sheet: widen column C, move selection ---
$wsCourse = $wb.Worksheets.Item("Course")
$wsCourse.Columns.Item(3).ColumnWidth = 22.15
[void]$wsCourse.Range("J8").Select()

# --- Course_Section sheet: add instructor_id column G ---
$wsCourseSection = $wb.Worksheets.Item("Course_Section")
$wsCourseSection.Columns.Item(7).ColumnWidth = 11.65
$wsCourseSection.Cells.Item(1, 7).Value = "instructor_id"
[void]$wsCourseSection.Cells.Item(1, 6).Copy()
[void]$wsCourseSection.Cells.Item(1, 7).PasteSpecial(-4122)
$wsCourseSection.Cells.Item(2, 7).Value = 10001
[void]$wsCourseSection.Range("G3").Select()

# --- Enrollment sheet: append row 2 with enrollment_id only ---
$wsEnrollment = $wb.Worksheets.Item("Enrollment")
$wsEnrollment.Cells.Item(2, 1).Value = 1
[void]$wsEnrollment.Range("A2").Select()

# --- Make Person the active sheet/tab (must be last so it "wins" as the active tab) ---
$wsPerson.Activate()
[void]$wsPerson.Range("B5").Select()
